$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "62.003.19"
$ws.Range("E2").Value = "  +4.71%  "

$ws.Range("D3").Value = "3.081.80"
$ws.Range("E3").Value = "  +2.75%  "

$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws.Range("D5") "580.37"
$ws.Range("E5").Value = "  +3.17%  "

Set-TextValue $ws.Range("D6") "142.32"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.070.84"
$ws.Range("E8").Value = "  +2.86%  "

$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("E10").Value = "  +4.81%  "

Set-TextValue $ws.Range("D11") "5.56"
$ws.Range("E11").Value = "  +8.31%  "

Set-TextValue $ws.Range("D12") "0.468"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("E13").Value = "  +4.12%  "

Set-TextValue $ws.Range("D14") "35.40"
$ws.Range("E14").Value = "  +4.97%  "

$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").Value = "3.588.18"
$ws.Range("E16").Value = "  +2.71%  "

Set-TextValue $ws.Range("D17") "7.30"
$ws.Range("E17").Value = "  +3.12%  "

$ws.Range("D18").Value = "3.077.19"
$ws.Range("E18").Value = "  +2.77%  "

$ws.Range("D19").Value = "61.905.49"
$ws.Range("E19").Value = "  +4.57%  "

Set-TextValue $ws.Range("D20") "449.43"
$ws.Range("E20").Value = "  +4.68%  "

Set-TextValue $ws.Range("D21") "13.93"
$ws.Range("E21").Value = "  +2.27%  "

Set-TextValue $ws.Range("D22") "0.730"
$ws.Range("E22").Value = "  +2.05%  "

Set-TextValue $ws.Range("D23") "7.44"
$ws.Range("E23").Value = "  +4.99%  "

$ws.Range("E24").Value = "  +2.96%  "

Set-TextValue $ws.Range("D25") "81.97"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("E27").Value = "  +5.23%  "

$ws.Range("E28").Value = "  -0.09%  "

Set-TextValue $ws.Range("D29") "2.67"
$ws.Range("E29").Value = "  +4.93%  "

Set-TextValue $ws.Range("D30") "8.16"
$ws.Range("E30").Value = "  +5.34%  "

Set-TextValue $ws.Range("D31") "6.77"
$ws.Range("E31").Value = "  +10.20%  "

Set-TextValue $ws.Range("D32") "0.112"
$ws.Range("E32").Value = "  +12.70%  "

Set-TextValue $ws.Range("D33") "26.79"
$ws.Range("E33").Value = "  +3.90%  "

$ws.Range("E34").Value = "  +4.34%  "

$ws.Range("D35").Value = "0.0₃0799"
$ws.Range("E35").Value = "  +2.78%  "

Set-TextValue $ws.Range("D36") "6.07"
$ws.Range("E36").Value = "  +3.72%  "

$ws.Range("E37").Value = "  +5.57%  "

Set-TextValue $ws.Range("D38") "50.39"
$ws.Range("E38").Value = "  +2.18%  "

Set-TextValue $ws.Range("D39") "3.00"
$ws.Range("E39").Value = "  +8.84%  "

Set-TextValue $ws.Range("D40") "8.83"
$ws.Range("E40").Value = "  +2.12%  "

Set-TextValue $ws.Range("D41") "430.18"
$ws.Range("E41").Value = "  +6.61%  "

$ws.Range("E42").Value = "  +5.58%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.834.35"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D44") "0.274"
$ws.Range("E44").Value = "  +7.75%  "

$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("E46").Value = "  +5.75%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D47") "0.999"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D48") "35.12"
$ws.Range("E48").Value = "  +4.47%  "

Set-TextValue $ws.Range("D49") "124.20"
$ws.Range("E49").Value = "  +1.11%  "

Set-TextValue $ws.Range("D50") "0.112"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("E51").Value = "  +2.48%  "
